$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("favorites")

# Delete rows 10 and 11 (they no longer exist in the updated table)
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()

# New data for rows 2-9 (A: media_type, B: media_id, C: favorite flag)
$data = @(
    @("movie", 496243, $false),
    @("movie", 278,    $false),
    @("tv",    1396,   $false),
    @("tv",    1429,   $false),
    @("movie", 496243, $true),
    @("movie", 278,    $true),
    @("tv",    1396,   $true),
    @("tv",    1429,   $true)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Update the selection to match the new state
$ws.Range("C6:C9").Select()
